# Update the line power values (columns B-F, rows 2-25) on Sheet1
# for the 380 kV case, per the recomputed results.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$data = @(
    @(1.301687528568038, 0.300747461727866, 0.1233630491873043, 0.4157543743867365, 2.798705207953788),
    @(1.171273146918452, 0.2623057578358043, 0.1128718144349108, 0.362076079561902, 2.580460778543483),
    @(1.092247776592785, 0.2388677213562289, 0.1064563759749859, 0.3292994308155812, 2.447576670238021),
    @(1.060301755946682, 0.229355877005446, 0.1038477432212375, 0.3159841952222422, 2.393692262385002),
    @(1.055012507639333, 0.2277787548360664, 0.1034149009966683, 0.3137755948689431, 2.384760503242802),
    @(1.091815907367675, 0.2387392848912668, 0.1064211730808466, 0.3291196943742705, 2.446848904307927),
    @(1.256498824658706, 0.2874570817163544, 0.1197396920338889, 0.3972057125687911, 2.723213523426864),
    @(1.588085390464755, 0.3844115841467897, 0.1461065271936377, 0.5323678854346667, 3.274723129474808),
    @(1.837444416842061, 0.456667782732552, 0.1656892636845839, 0.6329670626702324, 3.686736685555786),
    @(1.952238773908789, 0.4897972460817073, 0.1746568592192546, 0.6790791895586068, 3.875874300919293),
    @(1.995912514945587, 0.5023827667379237, 0.1780622397111813, 0.6965962014160993, 3.947759482263052),
    @(1.986497406172077, 0.4996704208222127, 0.1773283890871369, 0.6928210541511817, 3.932265758005713),
    @(1.955827710953031, 0.4908318444166184, 0.1749368246897234, 0.6805191853053856, 3.881782971614427),
    @(1.937068397614212, 0.4854232683167652, 0.1734731962532408, 0.6729912991885101, 3.850895548665449),
    @(1.829970381361761, 0.454508155220708, 0.1651044902739898, 0.6299609881346697, 3.674411954279094),
    @(1.764623137037574, 0.4356112001661359, 0.1599864691876292, 0.6036563113331539, 3.566595007603013),
    @(1.727164610783916, 0.4247664505884927, 0.1570482393692032, 0.5885590947900852, 3.504741527073719),
    @(1.714503450146765, 0.4210986920562618, 0.1560543227081723, 0.5834528703535824, 3.483825968738444),
    @(1.771566211276365, 0.4376202807796403, 0.1605307126019682, 0.6064530846480238, 3.578055627605522),
    @(1.964830552200908, 0.4934268375054671, 0.175639017821851, 0.6841309983901596, 3.896603718956044),
    @(2.092330707953749, 0.5301347953886193, 0.1855693007810544, 0.7352231267241507, 4.106333638094725),
    @(2.024169780583463, 0.5105206166376774, 0.1802638442694331, 0.7079228058401696, 3.99425021284344),
    @(1.768426903915952, 0.436711915334115, 0.1602846471847954, 0.6051885837940887, 3.572873871679349),
    @(1.497409031274401, 0.3580147870183055, 0.1389408687351619, 0.4955969926831187, 3.124400418051067)
)

$startRow = 2
for ($i = 0; $i -lt $data.Length; $i++) {
    $rowValues = $data[$i]
    $rowIndex = $startRow + $i
    for ($j = 0; $j -lt $rowValues.Length; $j++) {
        $colIndex = 2 + $j  # Column B = 2
        $ws.Cells.Item($rowIndex, $colIndex).Value = $rowValues[$j]
    }
}
